# Update the "Instructions" sheet of the FWC Calculator data template to
# reflect the revised app UI copy: renumbered/reworded steps, a new
# "Reload App" step, an updated sheet-count requirement (now three sheets,
# including the Instructions sheet itself), and the min/max composite
# volume note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# --- Step 1 text tweak: drop the word "here" ---
$ws.Range("A3").Value = "1. Download the Excel template file and overwrite it with your data. See the Data Requirements section below."

# --- Step 3 text tweak: "Start/End Time" -> "Start/End Date/Time", drop the elapsed-time sentence ---
$ws.Range("A7").Value = "3. Use the 'Start Date/Time' and 'End Date/Time' inputs to filter the data to the appropriate time range. The grayed-out sections of the graph will not be included in the aliquot volume and event mean concentration calculations."

# --- Step 4 text tweak: "Start/End Time" -> "Start/End Date/Time", "Redraw Graph(s)" -> "Draw Graph(s)" ---
$ws.Range("A8").Value = "4. After changing the 'Start Date/Time' and 'End Date/Time' inputs, click the 'Draw Graph(s)' button to regenerate the aliquot volume table, hydrograph, and pollutograph(s), filtered to the provided times."

# --- Step 5 text tweak: add min/max composite volume note ---
$ws.Range("A9").Value = "5. The 'Composite Vol.' input is used in the aliquot volume calculation such that the sum of the aliquot volumes will be equal to the composite volume value entered here, measured in mL. The minimum and maximum supported values are 500 mL and 10,000 mL, respectively."

# --- Insert new Step 7 ("Reload App") right after the existing Step 6 (row 10) ---
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "7. Use the 'Reload App' button to submit a new data set."

# --- Data Requirements: now three sheets instead of two (row 14 -> row 15 after the insert above) ---
$ws.Range("A15").Value = "  * Must contain exactly three sheets, in the following order:"

# --- Insert a new bullet describing the Instructions sheet itself, right after the "three sheets" bullet ---
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "    - Instructions: instructions for using the calculator"
